$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: new columns I ("I0") and J ("IF") ---
# Copy H1 formatting (bold, centered, thin border) onto I1:J1 first,
# then overwrite the text so the new headers match the look of the others.
$ws.Range("H1").Copy($ws.Range("I1:J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows 2-35: (row, I0 value, IF value) ---
$data = @(
  @(2, 7, 7),
  @(3, 8, 8),
  @(4, 1, 5),
  @(5, 1, 5),
  @(6, 1, 5),
  @(7, 1, 6),
  @(8, 1, 3),
  @(9, 1, 6),
  @(10, 1, 7),
  @(11, 1, 5),
  @(12, 1, 4),
  @(13, 1, 5),
  @(14, 1, 6),
  @(15, 1, 5),
  @(16, 1, 6),
  @(17, 1, 2),
  @(18, 1, 5),
  @(19, 1, 5),
  @(20, 1, 4),
  @(21, 1, 7),
  @(22, 1, 5),
  @(23, 1, 7),
  @(24, 1, 5),
  @(25, 1, 5),
  @(26, 1, 6),
  @(27, 1, 5),
  @(28, 1, 5),
  @(29, 1, 5),
  @(30, 1, 2),
  @(31, 7, 8),
  @(32, 5, 7),
  @(33, 4, 5),
  @(34, 1, 1),
  @(35, 1, 1)
)

foreach ($row in $data) {
  $r = $row[0]
  $ws.Cells.Item($r, 9).Value = $row[1]   # column I = I0
  $ws.Cells.Item($r, 10).Value = $row[2]  # column J = IF
}
